$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: was branch_and_bound('Vicky','JoseCid',L) -> now Joao/JoseCid test ---
$ws.Range("C20").Value = "branch_and_bound('Joao','JoseCid',L)"
$ws.Range("D20").Value = "L=[' Joao','Tiago','Stephanie','Diogo','Francisco','JoseCid']                                  "
$ws.Range("E20").Value = "L=[' Joao','Tiago','Stephanie','Diogo','Francisco','JoseCid']    "

# --- Row 21: branch_and_bound('Simao','Maria',L) stays, but D/E results change ---
$ws.Range("C21").Value = "branch_and_bound('Simao','Maria',L)"
$ws.Range("D21").Value = "L=[Simao','Joao','Tiago','Stephanie','Maria']"
$ws.Range("E21").Value = "L=[Simao','Artur','Tiago','Stephanie','Maria']"

# --- Row 22: branch_and_bound('Joao','Tiago',L) stays, but D/E result changes ---
$ws.Range("C22").Value = "branch_and_bound('Joao','Tiago',L)"
$ws.Range("D22").Value = "L=['Joao','Tiago']"
$ws.Range("E22").Value = "L=['Joao','Tiago']"

# --- Row 23 (new): branch_and_bound('Tiago','Joao',L) ---
$ws.Range("B23").Value = "Knowledge basis"
$ws.Range("C23").Value = "branch_and_bound('Tiago','Joao',L)"
$ws.Range("D23").Value = "L=['Tiago','Joao']"

# --- Row 24 (new): branch_and_bound('Andre','Nando') ---
$ws.Range("B24").Value = "Knowledge basis"
$ws.Range("C24").Value = "branch_and_bound('Andre','Nando')"
$ws.Range("D24").Value = "L=['Andre','Tiago','Stephanie','Diogo','Joao','Simao','Artur','Alcides','Nando']"
$ws.Range("E24").Value = "L=['Andre','Tiago','Stephanie','Diogo','Joao','Simao','Artur','Alcides','Nando']"

# --- back to row 23's E cell (matches original authoring order -> shared-string table order) ---
$ws.Range("E23").Value = "L=['Tiago','Stephanie','Diogo','Joao']"

# --- Row 25 (new): branch_and_bound('JoseCid','Maria') ---
$ws.Range("B25").Value = "Knowledge basis"
$ws.Range("C25").Value = "branch_and_bound('JoseCid','Maria')"
$ws.Range("D25").Value = "L=['JoseCid','Francisco','Diogo','Stephanie','Maria']"
$ws.Range("E25").Value = "L=['JoseCid','Francisco','Diogo','Tiago','Stephanie','Maria']"

# --- Row 30: note moved up (string re-shuffled in shared strings, text unchanged) ---
$ws.Range("B30").Value = "Note: Some bugs of branch_and_bound need to be fixed"

# --- Column C width widened to fit new, longer text ---
$ws.Range("C1").EntireColumn.ColumnWidth = 46.6

# --- Sheet view: drop the scrolled topLeftCell and move the selection to E25 ---
$ws.Range("E25").Select()
